$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (column D) values per the symbol-list refresh commit.
# Cells are stored as text, so force the Text number format before
# writing the value -- this preserves exact decimal formatting
# (e.g. trailing zeros like "22.40") instead of Excel coercing the
# literal into a numeric value.
$priceUpdates = @{
    "D2" = "248.23"
    "D3" = "22.40"
    "D4" = "5.363"
    "D5" = "0.05703"
    "D6" = "3.396"
    "D7" = "6.317"
    "D9" = "0.9130"
    "D10" = "0.1412"
    "D11" = "0.07395"
    "D12" = "0.03093"
    "D13" = "0.03017"
    "D14" = "0.09355"
    "D15" = "3.721"
    "D16" = "0.001578"
    "D17" = "0.04760"
    "D18" = "0.01829"
    "D19" = "0.0005792"
    "D20" = "0.006497"
    "D21" = "0.005004"
    "D23" = "0.0001500"
    "D24" = "3.697"
    "D25" = "2.154"
    "D26" = "0.3312"
    "D40" = "0.03972"
    "D41" = "0.006873"
    "D43" = "0.002711"
    "D44" = "0.007450"
    "D45" = "0.00005892"
    "D46" = "0.00000000750"
    "D47" = "0.5002"
    "D49" = "0.00002101"
    "D50" = "0.01010"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$cellRef]
}

Write-Host "Updated $($priceUpdates.Count) price cells"
